$d = $word.ActiveDocument

# Locate the sentence that mentions the number of years studied at UNO.
$sentence = $d.Range($d.Content.Start, $d.Content.End)
$sentence.Find.Execute("I studied at UNO for 2 years")

# Narrow down to just the digit(s) that need to change ("2" -> "10").
$digits = $d.Range($sentence.Start, $sentence.End)
$digits.Find.Execute("2")
$digitsStart = $digits.Start
$digitsEnd = $digits.End

# Replace "2" with "10" in place.
$numRange = $d.Range($digitsStart, $digitsEnd)
$numRange.Text = "10"
$numEnd = $digitsStart + 2   # end of the newly-typed "10"

# Force a run boundary right before the new text using a throwaway bookmark,
# so "I studied at UNO for " and "10" end up as separate runs.
$preSplit = $d.Range($digitsStart, $digitsStart)
$d.Bookmarks.Add("zzz_tmp_split", $preSplit)

# Drop Word's "_GoBack" bookmark exactly where the edit ended (right after
# the newly typed "10", before " years") - mirrors what Word does after a
# real edit.
$goBack = $d.Range($numEnd, $numEnd)
$d.Bookmarks.Add("_GoBack", $goBack)

# Remove the throwaway bookmark now that the run split exists; the real
# paragraph that used to hold "_GoBack" is left empty.
$d.Bookmarks("zzz_tmp_split").Delete()
